$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 - this shifts rows 3..19 up to 2..18, naturally dropping the
# old last row (19) and updating A/B/C/D columns (and any stale E values)
# to match the new "one year earlier" window the forecaster now reports.
$ws.Rows("2:2").Delete()

# The E column ("y_1_forecast") is recomputed by the naive forecaster and
# does not simply shift with the other columns, so refresh it explicitly
# for every data row, clearing rows that should now be blank.
$eValues = @{
    2 = $null
    3 = $null
    4 = $null
    5 = $null
    6 = 1.490702606731831
    7 = 3.498411812952029
    8 = 4.543069198269034
    9 = 4.42512316868644
    10 = 4.210645455310114
    11 = 4.757571096183799
    12 = 4.569144243718659
    13 = 3.124801698476176
    14 = -4.101394328717845
    15 = 3.801772939051373
    16 = 0.9099262091262217
    17 = 0.1232424362653362
    18 = 1.743978804508384
}

foreach ($r in $eValues.Keys) {
    $v = $eValues[$r]
    if ($null -eq $v) {
        $ws.Cells.Item($r, "E").ClearContents()
    } else {
        $ws.Cells.Item($r, "E").Value = $v
    }
}

# A couple of the "C" (y_0_forecast) values were recomputed with slightly
# different floating point precision; make sure they match exactly.
$ws.Cells.Item(4, "C").Value = 6.130685532900881
$ws.Cells.Item(5, "C").Value = 8.703939237319025

# The table now ends at row 18 instead of 19; make sure nothing lingers
# in row 19 so the sheet dimension shrinks to the smaller range.
$ws.Rows("19:19").ClearContents()
